$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlPasteValues  = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper: assign text that "looks numeric" (e.g. a percentage string) to a
# cell without letting the engine auto-convert it into a real number. We do
# this by writing the text as a literal-string formula and then collapsing
# the formula down to its cached value via PasteSpecial (values only), which
# keeps the existing cell style untouched and leaves a plain text cell.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial($xlPasteValues) | Out-Null
    $excel.CutCopyMode = 0
}

# --- G2: reorder the "Recorded By" email list (same people, new order) ---
$ws.Range("G2").Value2 = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

# --- Row 3: the ANATOMY / session 2 row moved from Pending -> Recorded ---
# Re-use the "Recorded" row look (style of row 2) for A3:F3 via a
# formats-only paste so the existing text values are left untouched.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A3:F3").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# Fill in who recorded it, the new attendance count and the new status
$ws.Range("G3").Value2 = "majorelle.magdy@med.asu.edu.eg"
$ws.Range("H3").Value2 = "19/251"
$ws.Range("I3").Value2 = "Recorded"

# Match the "Recorded" formatting on G3:I3 as well (copied from G2:I2)
$ws.Range("G2:I2").Copy() | Out-Null
$ws.Range("G3:I3").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Class Statistics block ---
$ws.Range("L6").Value2 = 6          # Recorded Sessions
$ws.Range("L8").Value2 = 23         # Pending Sessions
Set-TextValue $ws.Range("L9") "20.7%"   # Coverage %
Set-TextValue $ws.Range("L10") "23.3%"  # Average Attendance %

# --- Group Statistics summary row (row 15) ---
$ws.Range("O15").Value2 = 6
$ws.Range("Q15").Value2 = 23
Set-TextValue $ws.Range("R15") "20.7%"
Set-TextValue $ws.Range("S15") "23.3%"
